# paises.xlsx -- "Pais" sheet COVID-19 snapshot refresh.
#
# The source dashboard was re-pulled about an hour later (11:04 -> 12:04), which
# brings new totals for "Casos totales" (B), "Nuevos casos" (C), "Casos activos" (D),
# "Recuperados" (E), "Casos criticos" (F), "Muertes hoy" (G) and "Muertes" (H) for a
# couple dozen countries. Column A (the country name) never moves -- each row keeps
# referring to the same country both before and after the refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "last updated" timestamp shown in the title row.
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 12:04"

# Row 18: Belgica
$row = New-Object 'object[,]' 1,7
$row[0,0] = 51420
$row[0,1] = 639
$row[0,2] = 12980
$row[0,3] = 30025
$row[0,4] = 538
$row[0,5] = 76
$row[0,6] = 8415
$ws.Range("B18:H18").Value = $row

# Row 37: Rumania
$row = New-Object 'object[,]' 1,7
$row[0,0] = 14499
$row[0,1] = 392
$row[0,2] = 6144
$row[0,3] = 7479
$row[0,4] = 234
$row[0,5] = 12
$row[0,6] = 876
$ws.Range("B37:H37").Value = $row

# Row 39: Indonesia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 12776
$row[0,1] = 338
$row[0,2] = 2381
$row[0,3] = 9465
$row[0,4] = 0
$row[0,5] = 35
$row[0,6] = 930
$ws.Range("B39:H39").Value = $row

# Row 60: Moldavia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 4476
$row[0,1] = 0
$row[0,2] = 1747
$row[0,3] = 2584
$row[0,4] = 237
$row[0,5] = 2
$row[0,6] = 145
$ws.Range("B60:H60").Value = $row

# Row 72: Uzbekistan
$row = New-Object 'object[,]' 1,7
$row[0,0] = 2267
$row[0,1] = 2
$row[0,2] = 1002
$row[0,3] = 1157
$row[0,4] = 12
$row[0,5] = 0
$row[0,6] = 108
$ws.Range("B72:H72").Value = $row

# Row 73: Camerun
$row = New-Object 'object[,]' 1,7
$row[0,0] = 2266
$row[0,1] = 33
$row[0,2] = 1592
$row[0,3] = 664
$row[0,4] = 8
$row[0,5] = 0
$row[0,6] = 10
$ws.Range("B73:H73").Value = $row

# Row 87: Eslovenia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 1449
$row[0,1] = 1
$row[0,2] = 247
$row[0,3] = 1103
$row[0,4] = 13
$row[0,5] = 0
$row[0,6] = 99
$ws.Range("B87:H87").Value = $row

# Row 92: Hong Kong
$row = New-Object 'object[,]' 1,7
$row[0,0] = 1045
$row[0,1] = 4
$row[0,2] = 944
$row[0,3] = 97
$row[0,4] = 1
$row[0,5] = 0
$row[0,6] = 4
$ws.Range("B92:H92").Value = $row

# Row 138: Islas Feroe
$row = New-Object 'object[,]' 1,7
$row[0,0] = 187
$row[0,1] = 25
$row[0,2] = 93
$row[0,3] = 90
$row[0,4] = 1
$row[0,5] = 0
$row[0,6] = 4
$ws.Range("B138:H138").Value = $row

# Row 139: Martinica
$row = New-Object 'object[,]' 1,7
$row[0,0] = 187
$row[0,1] = 0
$row[0,2] = 185
$row[0,3] = 2
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$ws.Range("B139:H139").Value = $row

# Row 140: Liberia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 182
$row[0,1] = 0
$row[0,2] = 83
$row[0,3] = 85
$row[0,4] = 5
$row[0,5] = 0
$row[0,6] = 14
$ws.Range("B140:H140").Value = $row

# Row 141: Santo Tome y Principe
$row = New-Object 'object[,]' 1,7
$row[0,0] = 178
$row[0,1] = 0
$row[0,2] = 75
$row[0,3] = 83
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 20
$ws.Range("B141:H141").Value = $row

# Row 142: Republica del Chad
$row = New-Object 'object[,]' 1,7
$row[0,0] = 174
$row[0,1] = 0
$row[0,2] = 4
$row[0,3] = 167
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 3
$ws.Range("B142:H142").Value = $row

# Row 143: Birmania
$row = New-Object 'object[,]' 1,7
$row[0,0] = 170
$row[0,1] = 0
$row[0,2] = 43
$row[0,3] = 110
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 17
$ws.Range("B143:H143").Value = $row

# Row 144: Etiopia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 162
$row[0,1] = 1
$row[0,2] = 50
$row[0,3] = 106
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 6
$ws.Range("B144:H144").Value = $row

# Row 157: Haiti
$row = New-Object 'object[,]' 1,7
$row[0,0] = 108
$row[0,1] = 7
$row[0,2] = 15
$row[0,3] = 81
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 12
$ws.Range("B157:H157").Value = $row

# Row 191: Nueva Caledonia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 18
$row[0,1] = 0
$row[0,2] = 16
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 2
$ws.Range("B191:H191").Value = $row

# Row 192: Belice
$row = New-Object 'object[,]' 1,7
$row[0,0] = 18
$row[0,1] = 0
$row[0,2] = 18
$row[0,3] = 0
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$ws.Range("B192:H192").Value = $row

# Row 196: Namibia
$row = New-Object 'object[,]' 1,7
$row[0,0] = 16
$row[0,1] = 0
$row[0,2] = 9
$row[0,3] = 7
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$ws.Range("B196:H196").Value = $row

# Row 198: Dominica
$row = New-Object 'object[,]' 1,7
$row[0,0] = 16
$row[0,1] = 0
$row[0,2] = 13
$row[0,3] = 2
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 1
$ws.Range("B198:H198").Value = $row

# Row 199: Curazao
$row = New-Object 'object[,]' 1,7
$row[0,0] = 16
$row[0,1] = 0
$row[0,2] = 14
$row[0,3] = 2
$row[0,4] = 0
$row[0,5] = 0
$row[0,6] = 0
$ws.Range("B199:H199").Value = $row
